$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1832
$ws1.Range("F5").Value = 7
$ws1.Range("F7").Value = 182
$ws1.Range("F12").Value = 5262
$ws1.Range("F14").Value = 874
$ws1.Range("F16").Value = 2339
$ws1.Range("F18").Value = 41
$ws1.Range("F19").Value = 2183

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1832
$ws4.Range("F5").Value = 7
$ws4.Range("F7").Value = 182
$ws4.Range("F12").Value = 5262
$ws4.Range("F16").Value = 874
$ws4.Range("F18").Value = 2339
$ws4.Range("F21").Value = 41
$ws4.Range("F22").Value = 2183
